# Daily scrape update - 2025-08-10 03:44:12 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows (previously rows 4-6); only two data rows remain.
$ws.Rows("4:6").Delete()

# Row 2 — new opportunity data
# Leading apostrophe keeps the numeric-looking ID stored as text (matching
# the source data's inlineStr type) instead of Excel auto-converting it to
# a number; re-applying the Normal style afterwards drops the transient
# quote-prefix formatting bit so the cell is plain, unstyled text.
$ws.Range("A2").Value = "'1326674"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326674"
$ws.Range("C2").Value = "Foodprint"
$ws.Range("D2").Value = "Vanadzor, Armenia"
$ws.Range("E2").Value = "No"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").Value = "3 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Eco Art House"

# Row 3 — new opportunity data
$ws.Range("A3").Value = "'1323761"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1323761"
$ws.Range("C3").Value = "Sales representative"
$ws.Range("D3").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "4 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Fekretk"

# Column width adjustments (ColumnWidth units are 5/6 narrower than the
# stored OOXML "width" value, so subtract that fixed offset to land exactly
# on the target stored width).
$offset = 0.8333333333333334
$ws.Columns("C").ColumnWidth = 23 - $offset
$ws.Columns("D").ColumnWidth = 70 - $offset
$ws.Columns("F").ColumnWidth = 15 - $offset
$ws.Columns("G").ColumnWidth = 15 - $offset
$ws.Columns("H").ColumnWidth = 16 - $offset
